# Edit slide 21 ("Implementing Error Recovery (continued)") content placeholder:
#  1. Split "Only five methods throw a " into "Only three " + "methods throw a "
#     (changes "five" -> "three" and introduces a run break at that point).
#  2. Merge the two runs "                       " and "//    and " into a
#     single run "                       //    and ".
#  3. Remove the two paragraphs for parseIndexExpr()/parseFieldExpr() entirely.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(21)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Change 1: "Only five " -> "Only three " (splits run into two) ---
$first10 = $tr.Characters(1, 10)
if ($first10.Text -ne "Only five ") {
    throw "Unexpected text at start of placeholder: [$($first10.Text)]"
}
$first10.Text = "Only three "

# --- Change 2: merge "                       " + "//    and " runs ---
$full = $tr.Text
$needle = "                       //    and "
$idx = $full.IndexOf($needle)
if ($idx -lt 0) {
    throw "Could not locate spacer/comment run text to merge"
}
$mergeRange = $tr.Characters($idx + 1, $needle.Length)
$mergeRange.Text = $needle

# --- Change 3: delete the parseIndexExpr()/parseFieldExpr() paragraphs ---
$paraIndex = $tr.Paragraphs(4, 1)
$paraField = $tr.Paragraphs(5, 1)
if ($paraIndex.Text -notlike "parseIndexExpr*") {
    throw "Unexpected paragraph 4 text: [$($paraIndex.Text)]"
}
if ($paraField.Text -notlike "parseFieldExpr*") {
    throw "Unexpected paragraph 5 text: [$($paraField.Text)]"
}
$paraAdd = $tr.Paragraphs(6, 1)
$delStart = $paraIndex.Start
$delLen = $paraAdd.Start - $paraIndex.Start
$toDelete = $tr.Characters($delStart, $delLen)
$toDelete.Delete()
